$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "for a sudoku solving algorithm, using a DLS approach ends up",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "for a sudoku solving algorithm using a DLS approach, ends up", 2
)
